$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the bulk of row 3's data (CLS + Thu Thuat work in progress) ---
# Remove all values from row 3 (A3:CB3) first.
$ws.Range("A3:CB3").ClearContents()

# AL3 had fill formatting tied to its old boolean value; the row rebuild drops it
# entirely, so clear its formatting along with its content.
$ws.Range("AL3").Clear()

# AB3 should still be present as a bare placeholder cell (matching the rest of the
# column, e.g. AB4, AB5, ...). Copying formats from a neighboring placeholder cell
# forces Excel to keep an explicit but style-less <c r="AB3"/> entry.
$ws.Range("AB4").Copy()
$ws.Range("AB3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view: scroll back to the top-left and move the selection ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("BW8").Select()
